# Adds a new "LDKit" / "LDkit" results column (K,L) to both worksheets,
# mirroring the existing 4-column-group layout (All (x100) / Search by
# field (x100)) used for LD-flex, and extends the summary rows.

$wb = $excel.ActiveWorkbook

$kValues = @(30211,32276,29153,27197,27243,39635,46214,45354,48396,47155,33749,31541,31988,32092,28821,27174,27146,27158,27150,27022,26986,27302,27656,27812,28815,27048,27003,28675,48463,48650)
$lValues = @(26356,25427,41354,49674,46745,45535,45269,46616,48880,46749,29423,30701,31287,31662,30557,31238,31726,30593,29899,30103,27740,24731,24073,24244,25263,26107,24012,24139,24038,25018)

# ============================================================
# Sheet "Measurements"
# ============================================================
$ws = $wb.Worksheets.Item("Measurements")

# --- column width for the new column L (same as column H) ---
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- row 1: merged header "LDKit" above K:L, same look as I1:J1 ("LD-flex") ---
$ws.Range("I1:J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)
$ws.Range("K1").Value = "LDKit"
$ws.Range("K1:L1").Merge()

# --- row 2: sub headers ---
$ws.Range("I2:J2").Copy()
$ws.Range("K2:L2").PasteSpecial(-4122)
$ws.Range("K2").Value = "All (x100)"
$ws.Range("L2").Value = "Search by field (x100)"

# --- rows 3-32: raw measurements ---
$ws.Range("I3:J32").Copy()
$ws.Range("K3:L32").PasteSpecial(-4122)
for ($i = 0; $i -lt 30; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 11).Value = $kValues[$i]
    $ws.Cells.Item($r, 12).Value = $lValues[$i]
}

# --- row 34: Mean (average of rows 3-32), extend the shared range B34:J34 -> B34:L34 ---
$ws.Range("I34:J34").Copy()
$ws.Range("K34:L34").PasteSpecial(-4122)
$ws.Range("B34:L34").Formula = "=AVERAGE(B3:B32)"

# --- row 35: Normalised (ms) ---
$ws.Range("I35:J35").Copy()
$ws.Range("K35:L35").PasteSpecial(-4122)
$ws.Range("K35:L35").Formula = "=K34/100"

# --- row 37 label reindex happens automatically (string content unchanged) ---

# --- row 42: new "LDkit" summary line, mirroring rows 38/40/41 ---
$ws.Range("A41:C41").Copy()
$ws.Range("A42:C42").PasteSpecial(-4122)
$ws.Range("A42").Value = "LDkit"
$ws.Range("B42:C42").Formula = "=K35"

# --- merge for the new header cell ---
# (merge already performed above for K1:L1)

# ============================================================
# Sheet "Normalised measurements"
# ============================================================
$ws2 = $wb.Worksheets.Item("Normalised measurements")

# --- row 1: header "LDkit" above K1 only (NOT merged, matches source quirk) ---
$ws2.Range("I2").Copy()
$ws2.Range("K1").PasteSpecial(-4122)
$ws2.Range("K1").Value = "LDkit"

# --- row 2: sub headers ---
$ws2.Range("I2:J2").Copy()
$ws2.Range("K2:L2").PasteSpecial(-4122)
$ws2.Range("K2").Value = "All (x100)"
$ws2.Range("L2").Value = "Search by field (x100)"

# --- rows 3-32: normalised (divide the raw measurement by 100) ---
$ws2.Range("I3:J32").Copy()
$ws2.Range("K3:L32").PasteSpecial(-4122)
for ($i = 0; $i -lt 30; $i++) {
    $r = 3 + $i
    $ws2.Cells.Item($r, 11).Formula = "=Measurements!K$r/100"
    $ws2.Cells.Item($r, 12).Formula = "=Measurements!L$r/100"
}

# --- row 34: Mean, extend the shared range B34:J34 -> B34:L34 ---
$ws2.Range("I34:J34").Copy()
$ws2.Range("K34:L34").PasteSpecial(-4122)
$ws2.Range("B34:L34").Formula = "=AVERAGE(B3:B32)"

# --- row 42: new "LDkit" summary line, mirroring rows 38/40/41 ---
$ws2.Range("A41:C41").Copy()
$ws2.Range("A42:C42").PasteSpecial(-4122)
$ws2.Range("A42").Value = "LDkit"
$ws2.Range("B42:C42").Formula = "=K34"

Write-Output "edit applied"
